$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.424.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.924.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.32%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.58"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.14"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.42"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.18"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.386.55"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.937.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.78%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.486.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.74"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.66"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.12"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.21"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "53.12"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.18%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0455"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.61%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.94"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.47%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +10.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.41"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +11.52%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.36%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.223.55"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.264"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +24.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0342"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +16.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.965"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.21%  "
